$wb = $excel.ActiveWorkbook

# --- Sheet "NBR" (first sheet) - update column C values (Reaction_number) ---
$wsNBR = $wb.Worksheets.Item("NBR")
$nbrValues = @(774,772,731,727,726,721,716,713,710,705,706,701,702,670,667,663,659,652,638)
for ($i = 0; $i -lt $nbrValues.Length; $i++) {
    $row = $i + 2
    $wsNBR.Cells.Item($row, 3).Value = $nbrValues[$i]
}

# --- Sheet "BAR" (second sheet) - update column C values (Reaction_number) ---
$wsBAR = $wb.Worksheets.Item("BAR")
$barValues = @(688,691,691,691,690,691,691,691,682,679,675,676,676,665,663,663,663,663,663)
for ($i = 0; $i -lt $barValues.Length; $i++) {
    $row = $i + 2
    $wsBAR.Cells.Item($row, 3).Value = $barValues[$i]
}
